$wb = $excel.ActiveWorkbook

$newMdName = "27144277-3f91-4e5e-8751-84420dbded78.md"
$newStatus = "Handoff transform failed"
$zeroDate  = "0001-01-01 00:00:00"
$ignored   = "Ignored"

# NOTE: the underlying hyperlink relationship target keeps pointing at the
# old markdown file (only the cell text / display caption changes to the
# new handoff file name) - matches the unchanged *.xml.rels in the diff.
$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/d7a17e171e413dd065a643b10b0596ee1f0ef1d0/e2e/f82fc655-b24f-4177-8998-9e0c2c61e621.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f714b04fcc4ba169d2dea18196c38d691af93a15/.localization-config"

# ---------------------------------------------------------------------------
# Overview sheet: new handoff markdown file name, new status text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet: handoff attempt failed/ignored, target+handoff file cleared.
# ---------------------------------------------------------------------------
$xlfUrlZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a50327f3b02931c00e82d0425b91714533f4cfc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/f82fc655-b24f-4177-8998-9e0c2c61e621.af4278979307631c9f7905d22e8f6148cdd6307a.zh-cn.xlf"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Clear()
$wsZh.Range("D2").Value = $zeroDate
$wsZh.Range("G2").Value = $zeroDate
$wsZh.Range("H2").Value = $ignored
$wsZh.Range("D3").Value = $zeroDate
$wsZh.Range("G3").Value = $zeroDate
$wsZh.Range("H3").Value = $ignored

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: same handling as zh-cn, different xlf locale.
# ---------------------------------------------------------------------------
$xlfUrlDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/492bc298ad85f860db2aeab4779127c37dd57205/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/f82fc655-b24f-4177-8998-9e0c2c61e621.af4278979307631c9f7905d22e8f6148cdd6307a.de-de.xlf"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Clear()
$wsDe.Range("D2").Value = $zeroDate
$wsDe.Range("G2").Value = $zeroDate
$wsDe.Range("H2").Value = $ignored
$wsDe.Range("D3").Value = $zeroDate
$wsDe.Range("G3").Value = $zeroDate
$wsDe.Range("H3").Value = $ignored

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null
